# Apply "Mastercard Orange" heading template update:
# - Heading1-4 styles get Calibri Light ascii/hAnsi fonts
# - Heading1-4 styles get solid color FF5F00 (replacing themed blue colors)
# - Heading1 grows from 14pt to 18pt; Heading2 from 13pt to 15pt
# - Heading3 and Heading4 gain an explicit font size (13pt / 12pt) that
#   they previously lacked (inherited from Normal)

$d = $word.ActiveDocument

# Mastercard Orange (RGB 255, 95, 0 / hex FF5F00), packed as a Word "OLE"
# BGR color value (0x00BBGGRR) the way RGB() would produce it:
#   (0x00 << 16) | (0x5F << 8) | 0xFF = 24575
$orange = 24575

$h1 = $d.Styles("Heading1")
$h1.Font.Name = "Calibri Light"
$h1.Font.Color = $orange
$h1.Font.Size = 18

$h2 = $d.Styles("Heading2")
$h2.Font.Name = "Calibri Light"
$h2.Font.Color = $orange
$h2.Font.Size = 15

$h3 = $d.Styles("Heading3")
$h3.Font.Name = "Calibri Light"
$h3.Font.Color = $orange
$h3.Font.Size = 13

$h4 = $d.Styles("Heading4")
$h4.Font.Name = "Calibri Light"
$h4.Font.Color = $orange
$h4.Font.Size = 12
